# Update "想去人数" (F column) figures across the four sheets to match
# the newly generated output (gh-pages data refresh at commit 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 10835
$ws.Range("F4").Value = 261
$ws.Range("F5").Value = 1198
$ws.Range("F6").Value = 1071
$ws.Range("F7").Value = 838
$ws.Range("F8").Value = 281
$ws.Range("F9").Value = 1149
$ws.Range("F12").Value = 880
$ws.Range("F13").Value = 381
$ws.Range("F14").Value = 1940
$ws.Range("F16").Value = 956
$ws.Range("F17").Value = 821
$ws.Range("F18").Value = 547
$ws.Range("F19").Value = 796
$ws.Range("F20").Value = 906
$ws.Range("F24").Value = 623
$ws.Range("F25").Value = 635
$ws.Range("F26").Value = 118
$ws.Range("F27").Value = 342
$ws.Range("F28").Value = 1014
$ws.Range("F31").Value = 167
$ws.Range("F34").Value = 564
$ws.Range("F35").Value = 1793
$ws.Range("F36").Value = 380
$ws.Range("F37").Value = 23
$ws.Range("F38").Value = 1424
$ws.Range("F42").Value = 84
$ws.Range("F44").Value = 2
$ws.Range("F45").Value = 76
$ws.Range("F46").Value = 82
$ws.Range("F47").Value = 40
$ws.Range("F48").Value = 8
$ws.Range("F49").Value = 81

# --- Sheet "演出" ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F14").Value = 136

# --- Sheet "本地生活" ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 2167
$ws.Range("F3").Value = 630
$ws.Range("F4").Value = 559

# --- Sheet "全部类型" ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 2167
$ws.Range("F3").Value = 630
$ws.Range("F5").Value = 10835
$ws.Range("F6").Value = 261
$ws.Range("F8").Value = 559
$ws.Range("F9").Value = 1071
$ws.Range("F11").Value = 1149
$ws.Range("F14").Value = 880
$ws.Range("F15").Value = 381
$ws.Range("F16").Value = 1940
$ws.Range("F18").Value = 956
$ws.Range("F19").Value = 821
$ws.Range("F20").Value = 547
$ws.Range("F21").Value = 796
$ws.Range("F22").Value = 906
$ws.Range("F27").Value = 623
$ws.Range("F30").Value = 635
$ws.Range("F31").Value = 118
$ws.Range("F32").Value = 342
$ws.Range("F33").Value = 1014
$ws.Range("F36").Value = 167
$ws.Range("F39").Value = 23
$ws.Range("F40").Value = 1424
$ws.Range("F44").Value = 84
$ws.Range("F46").Value = 76
$ws.Range("F47").Value = 40
$ws.Range("F48").Value = 81
